$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'-3.45%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'40.20"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-4.16%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.046"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-2.66%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07594"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-6.47%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'4.261"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-2.40%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'1.595"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-9.49%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.9078"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-2.55%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.1001"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-10.68%"
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'-5.42%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.09042"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-1.84%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.04373"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-6.10%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.1054"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'0.12%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001233"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-4.06%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.005863"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'2.53%"
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'2,406.14%"
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = "'0.31%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.459"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-3.61%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.3281"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-3.13%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'6.881"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-7.03%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.1361"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'-2.32%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.2825"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'8.52%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04164"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-0.72%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001212"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-2.66%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.004059"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-4.56%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'6.31%"
$ws.Range("E26").ClearFormats()
$ws.Range("D38").Value = "'0.02415"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'-6.60%"
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.05139"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'-6.33%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.007850"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'-3.50%"
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'-6.26%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.007088"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'8.52%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.001973"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-5.87%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.008378"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'9.60%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.3315"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-4.28%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006464"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-4.73%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'-0.02%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.004828"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'42.65%"
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'-27.00%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.00002109"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.0002009"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'-0.02%"
$ws.Range("E51").ClearFormats()
